$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the numeric Price/Volume columns (D:E) so values round-trip
# as literal text (matching the source data, which stores numbers/percentages as strings)
# rather than being auto-coerced into numeric/percentage cell types.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "316.03"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.67%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.38"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.28%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.124"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.44%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08174"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.88%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.982"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.66%"

$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.366"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.36%"

$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.331"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.09%"

$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9388"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.08%"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1307"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-7.52%"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1973"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.97%"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09053"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.03%"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03498"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.18%"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09738"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.80%"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001409"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.69%"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005974"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.94%"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.650"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-7.75%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.85%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.85%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.30%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.958"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "6.74%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2491"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.73%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04368"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.55%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001242"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.99%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004765"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "9.87%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003892"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "199.04%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-7.63%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02208"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "8.89%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05193"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.92%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007760"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.11%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01037"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "6.04%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1402"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "3.02%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002101"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.48%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009275"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "6.43%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006942"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "9.07%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.00%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002883"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.58%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001692"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "30.07%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.00%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.00%"
